# Add RADAR4Chem and Chemotion Repository
# Inserts two new rows (new rows 54 & 55) into the "Repositories" worksheet,
# just above the existing "Crystallography Open Database" entries, shifting
# all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repositories")
$ws.Activate()

# Insert two blank rows at row 54 (existing rows 54.. shift down to 56..)
$insertRange = $ws.Range("A54:J55")
$insertRange.EntireRow.Insert()

# Row 54: Chemotion Repository
$ws.Cells.Item(54, 1).Value = "KIT"
$ws.Cells.Item(54, 2).Value = "Chemotion Repository"
$ws.Cells.Item(54, 3).Value = "Repository"
$ws.Cells.Item(54, 4).Value = "Chemistry"
$ws.Cells.Item(54, 5).Value = "yes"
$ws.Cells.Item(54, 6).Value = "yes"
$ws.Cells.Item(54, 9).Value = "https://www.chemotion-repository.net"
$ws.Cells.Item(54, 10).Value = "Chemotion Repository provides an archive for materials, as well as research data repository for samples, reactions and analyses"

# Row 55: RADAR4Chem
$ws.Cells.Item(55, 1).Value = "FIZ Karlsruhe"
$ws.Cells.Item(55, 2).Value = "RADAR4Chem"
$ws.Cells.Item(55, 3).Value = "Repository"
$ws.Cells.Item(55, 4).Value = "Chemistry"
$ws.Cells.Item(55, 5).Value = "yes"
$ws.Cells.Item(55, 6).Value = "yes"
$ws.Cells.Item(55, 9).Value = "https://radar4chem.radar-service.eu"
$ws.Cells.Item(55, 10).Value = "API access only includes metadata"

# Reflect the author's final view state: scrolled down with rows 54:55 selected
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Range("A54:XFD55").Select()
